$wb = $excel.ActiveWorkbook

# The edited sheet is "Dataset_ex" (already the selected/active tab)
$ws = $wb.Worksheets.Item("Dataset_ex")
$ws.Activate()

# Update column A for rows 102..151: value becomes (row number - 1)
for ($r = 102; $r -le 151; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 1
}

# Update the view: scroll so that row 127 is the top-left visible row,
# and select cell C144
$appWin = $excel.ActiveWindow
$appWin.ScrollRow = 127
$appWin.ScrollColumn = 1
$ws.Range("C144").Select()
